$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (Q) mirroring the existing 2019 column (P):
# copy each source cell's formatting/value onto the new column, then
# overwrite the values that actually differ for 2020.
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 2

$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 0.3

$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("Q7").Value = 0.1

$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("Q8").Value = 4.3

# Match the recorded selection state left behind in the sheet view.
$ws.Activate()
$ws.Range("O12").Select()
